$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row labels: "_old" -> "_FV2310", "_new" -> "_FV2404" ---
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2310"
}

# Column 11 (K) is "diff" - left unchanged.

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2404"
}

# --- Turn the used range into an Excel Table named "Table1" ---
$tableRange = $ws.Range("A1:U87")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1, $null)
$lo.Name = "Table1"

# --- Freeze the header row (split after row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
